$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that looks like a plain number but must stay a text string
# (the source workbook stores every data cell as inline/shared text, never as a
# real number). Temporarily force the cell to Text format so Excel doesn't
# auto-convert the string to a numeric value, then strip the format change
# back off so no extra style gets attached to the cell.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# Delete row 7 entirely (the old "Сумма:" totals row) - this shrinks the
# used range from A1:K7 down to A1:K6 and shifts row 7's removal into place.
$ws.Rows.Item(7).Delete()

# ---- Row 2 ----
$ws.Range("B2").Value = "1.1 - 2.1"
Set-TextValue "C2" "1"
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""

# ---- Row 3 ----
$ws.Range("B3").Value = "28.1 - 31.5"
Set-TextValue "C3" "123"
$ws.Range("F3").Value = ""

# ---- Row 4 ----
$ws.Range("B4").Value = "27.7 - 27.7"
Set-TextValue "C4" "1"
$ws.Range("F4").Value = "asdfsa23"
$ws.Range("J4").Value = "2df"

# ---- Row 5 (cleared out) ----
$ws.Range("B5").Value = ""
$ws.Range("C5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("G5").Value = ""
$ws.Range("H5").Value = ""
$ws.Range("J5").Value = ""
$ws.Range("K5").Value = ""

# ---- Row 6 (now the totals "Сумма:" row) ----
$ws.Range("B6").Value = "Сумма:"
Set-TextValue "F6" "0"
Set-TextValue "G6" "234"
Set-TextValue "H6" "234"
Set-TextValue "K6" "234"
